$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws3 = $wb.Worksheets.Item("Tabelle3")

# --- Tabelle1: fill in the Netzplan (CPM) calculation cells with formulas ---
$ws = $ws1
$ws.Range("U2").Formula = "=R8"
$ws.Range("W2").Formula = "=U2+U4"
$ws.Range("V4").Formula = "=W5-W2"
$ws.Range("W4").Formula = "=AE8-W2"
$ws.Range("U5").Formula = "=W5-U4"
$ws.Range("W5").Formula = "=AE11"
$ws.Range("P8").Formula = "=M11"
$ws.Range("R8").Formula = "=P8+P10"
$ws.Range("U8").Formula = "=R8"
$ws.Range("W8").Formula = "=U8+U10"
$ws.Range("AE8").Formula = "=MAX(W2,W14,AB20)"
$ws.Range("AG8").Formula = "=AE8+AE10"
$ws.Range("Q10").Formula = "=R11-R8"
$ws.Range("R10").Formula = "=MIN(U2,U8,U14)-R8"
$ws.Range("V10").Formula = "=W11-W8"
$ws.Range("W10").Formula = "=Z20-W8"
$ws.Range("AF10").Formula = "=AG11-AG8"
$ws.Range("AG10").Value = 0
$ws.Range("M11").Formula = "=K11+K13"
$ws.Range("P11").Formula = "=R11-P10"
$ws.Range("R11").Formula = "=MIN(U5,U11,U17)"
$ws.Range("U11").Formula = "=W11-U10"
$ws.Range("W11").Formula = "=Z23"
$ws.Range("AE11").Formula = "=AG11-AE10"
$ws.Range("AG11").Formula = "=AG8"
$ws.Range("K14").Formula = "=M14-K13"
$ws.Range("M14").Formula = "=MIN(P11,P17)"
$ws.Range("P14").Formula = "=M11"
$ws.Range("R14").Formula = "=P14+P16"
$ws.Range("U14").Formula = "=MAX(R8,R14)"
$ws.Range("W14").Formula = "=U14+U16"
$ws.Range("Q16").Formula = "=R17-R14"
$ws.Range("R16").Formula = "=MIN(U14,Z20)-R14"
$ws.Range("V16").Formula = "=W17-W14"
$ws.Range("W16").Formula = "=AE8-W14"
$ws.Range("P17").Formula = "=R17-P16"
$ws.Range("R17").Formula = "=MIN(U17,Z23)"
$ws.Range("U17").Formula = "=W17-U16"
$ws.Range("W17").Formula = "=AE11"
$ws.Range("Z20").Formula = "=MAX(W8,R14)"
$ws.Range("AB20").Formula = "=Z20+Z22"
$ws.Range("AA22").Formula = "=AB23-AB20"
$ws.Range("AB22").Formula = "=AE8-AB20"
$ws.Range("Z23").Formula = "=AB23-Z22"
$ws.Range("AB23").Formula = "=AE11"

# --- selections per sheet (also drives which sheet tab is active) ---
$ws1.Range("R16").Select()
$ws3.Range("K11").Select()
$ws2.Range("V21").Select()
